$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")

# Row 80
$ws.Range("H80").Value = 208.4
$ws.Range("I80").Value = 208.4
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 625.2
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = 372.8
$ws.Range("N80").ClearContents()

# Row 83
$ws.Range("H83").Value = 208.4
$ws.Range("I83").Value = 208.4
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 1875.6
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = 3116.4
$ws.Range("N83").ClearContents()

# Row 100
$ws.Range("H100").Value = 12822975
$ws.Range("I100").Value = 20835460
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 20835460
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -20834919
$ws.Range("N100").Value = -4082

# Row 129
$ws.Range("H129").Value = 1056.1571
$ws.Range("I129").Value = 790
$ws.Range("J129").Value = 1076.6307
$ws.Range("K129").Value = 2370
$ws.Range("L129").Value = 3229.8921
$ws.Range("M129").Value = 2630
$ws.Range("N129").Value = -13229.8921

# Row 132
$ws.Range("H132").Value = 1168.7693
$ws.Range("I132").Value = 1022.1875
$ws.Range("J132").Value = 2927.75
$ws.Range("K132").Value = 3066.5625
$ws.Range("L132").Value = 8783.25
$ws.Range("M132").Value = -536.5625

# Row 137
$ws.Range("H137").Value = 2424.3667
$ws.Range("I137").Value = 2654.3333
$ws.Range("J137").Value = 2079.4167
$ws.Range("K137").Value = 7962.999899999999
$ws.Range("L137").Value = 6238.250100000001
$ws.Range("M137").Value = -5412.999899999999

# Row 138
$ws.Range("H138").Value = 1611.921
$ws.Range("I138").Value = 1332.4
$ws.Range("J138").Value = 1922.5
$ws.Range("K138").Value = 3997.2
$ws.Range("L138").Value = 5767.5
$ws.Range("M138").Value = 1142.8
$ws.Range("N138").Value = -16047.5


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")

# Row 2
$ws.Range("H2").Value = 3000
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 3000
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 3000
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -3226

# Row 9
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

# Row 20
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

# Row 45
$ws.Range("H45").Value = 7310.6875
$ws.Range("I45").Value = 8767
$ws.Range("J45").Value = 1000
$ws.Range("K45").Value = 8767
$ws.Range("L45").Value = 1000
$ws.Range("M45").Value = -8390
$ws.Range("N45").Value = -1754

# Row 74
$ws.Range("H74").Value = 2213.75
$ws.Range("I74").Value = 2065.375
$ws.Range("J74").Value = 2510.5
$ws.Range("K74").Value = 2065.375
$ws.Range("L74").Value = 2510.5
$ws.Range("M74").Value = -1191.375
$ws.Range("N74").Value = -4258.5

# Row 77
$ws.Range("H77").Value = 2213.75
$ws.Range("I77").Value = 2065.375
$ws.Range("J77").Value = 2510.5
$ws.Range("K77").Value = 10326.875
$ws.Range("L77").Value = 12552.5
$ws.Range("M77").Value = -5958.875
$ws.Range("N77").Value = -21288.5

# Row 116
$ws.Range("H116").Value = 3000
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 3000
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -7588


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")

# Row 3
$ws.Range("H3").Value = 3000
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 3000
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 3000
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -3228

# Row 94
$ws.Range("H94").Value = 1906.875
$ws.Range("I94").Value = 1414.7858
$ws.Range("J94").Value = 2595.8
$ws.Range("K94").Value = 1414.7858
$ws.Range("L94").Value = 2595.8
$ws.Range("M94").Value = -963.7858000000001
$ws.Range("N94").Value = -3497.8

# Row 105
$ws.Range("H105").Value = 2727.1667
$ws.Range("I105").Value = 1455.5454
$ws.Range("J105").Value = 3803.1538
$ws.Range("K105").Value = 1455.5454
$ws.Range("L105").Value = 3803.1538
$ws.Range("M105").Value = 291.4546
$ws.Range("N105").Value = -7297.1538

# Row 134
$ws.Range("H134").Value = 4609.3335
$ws.Range("I134").Value = 5386.2144
$ws.Range("J134").Value = 2631.818
$ws.Range("K134").Value = 16158.6432
$ws.Range("L134").Value = 7895.454000000001
$ws.Range("M134").Value = -13623.6432
$ws.Range("N134").Value = -12965.454


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")

# Row 25
$ws.Range("H25").Value = 142858000
$ws.Range("I25").Value = 1000
$ws.Range("J25").Value = 1000000000
$ws.Range("K25").Value = 1000
$ws.Range("L25").Value = 1000000000
$ws.Range("M25").Value = -826
$ws.Range("N25").Value = -1000000348


# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")

# Row 22
$ws.Range("H22").Value = 2522.4614
$ws.Range("I22").Value = 800
$ws.Range("J22").Value = 2835.6365
$ws.Range("K22").Value = 2400
$ws.Range("L22").Value = 8506.9095
$ws.Range("M22").Value = -2231
$ws.Range("N22").Value = -8844.9095

# Row 27
$ws.Range("H27").Value = 2522.4614
$ws.Range("I27").Value = 800
$ws.Range("J27").Value = 2835.6365
$ws.Range("K27").Value = 2400
$ws.Range("L27").Value = 8506.9095
$ws.Range("M27").Value = -2298
$ws.Range("N27").Value = -8710.9095

# Row 34
$ws.Range("H34").Value = 2006.3636
$ws.Range("I34").Value = 1014
$ws.Range("J34").Value = 2833.3333
$ws.Range("K34").Value = 3042
$ws.Range("L34").Value = 8499.999899999999
$ws.Range("M34").Value = -2958
$ws.Range("N34").Value = -8667.999899999999


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")

# Row 102
$ws.Range("H102").Value = 1131350.1
$ws.Range("I102").Value = 1884116.9
$ws.Range("J102").Value = 2200
$ws.Range("K102").Value = 1884116.9
$ws.Range("L102").Value = 2200
$ws.Range("M102").Value = -1882494.9

# Row 107
$ws.Range("H107").Value = 1741.4783
$ws.Range("I107").Value = 617.53845
$ws.Range("J107").Value = 3202.6
$ws.Range("K107").Value = 617.53845
$ws.Range("L107").Value = 3202.6
$ws.Range("M107").Value = 1302.46155
$ws.Range("N107").Value = -7042.6

# Row 122
$ws.Range("H122").Value = 3603867
$ws.Range("I122").Value = 10804636
$ws.Range("J122").Value = 3482.3333
$ws.Range("K122").Value = 32413908
$ws.Range("L122").Value = 10446.9999
$ws.Range("M122").Value = -32411458
$ws.Range("N122").Value = -15346.9999

# Row 126
$ws.Range("H126").Value = 4866.2925
$ws.Range("I126").Value = 9377
$ws.Range("J126").Value = 2772.0356
$ws.Range("K126").Value = 28131
$ws.Range("L126").Value = 8316.106800000001
$ws.Range("M126").Value = -25661
$ws.Range("N126").Value = -13256.1068


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")

# Row 93
$ws.Range("H93").Value = 1277.6666
$ws.Range("I93").Value = 999.5
$ws.Range("J93").Value = 1357.1428
$ws.Range("K93").Value = 999.5
$ws.Range("L93").Value = 1357.1428
$ws.Range("M93").Value = 248.5
$ws.Range("N93").Value = -3853.1428

# Row 100
$ws.Range("H100").Value = 1416.6666
$ws.Range("I100").Value = 934
$ws.Range("J100").Value = 2020
$ws.Range("K100").Value = 934
$ws.Range("L100").Value = 2020
$ws.Range("M100").Value = -393
$ws.Range("N100").Value = -3102

# Row 132
$ws.Range("H132").Value = 11910104
$ws.Range("I132").Value = 16672546
$ws.Range("J132").Value = 4000.375
$ws.Range("K132").Value = 50017638
$ws.Range("L132").Value = 12001.125
$ws.Range("M132").Value = -50015108
$ws.Range("N132").Value = -17061.125

# Row 136
$ws.Range("H136").Value = 6896.1353
$ws.Range("I136").Value = 5690.067
$ws.Range("J136").Value = 12065
$ws.Range("K136").Value = 17070.201
$ws.Range("L136").Value = 36195
$ws.Range("M136").Value = -14520.201
$ws.Range("N136").Value = -41295


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")

# Row 81
$ws.Range("H81").Value = 1364.1
$ws.Range("I81").Value = 1293.6666
$ws.Range("J81").Value = 1998
$ws.Range("K81").Value = 2587.3332
$ws.Range("L81").Value = 3996
$ws.Range("M81").Value = -1526.3332
$ws.Range("N81").Value = -6118

# Row 84
$ws.Range("H84").Value = 1364.1
$ws.Range("I84").Value = 1293.6666
$ws.Range("J84").Value = 1998
$ws.Range("K84").Value = 12936.666
$ws.Range("L84").Value = 19980
$ws.Range("M84").Value = -7632.666000000001
$ws.Range("N84").Value = -30588

# Row 100
$ws.Range("H100").Value = 72131.64
$ws.Range("I100").Value = 167192.17
$ws.Range("J100").Value = 836.25
$ws.Range("K100").Value = 334384.34
$ws.Range("L100").Value = 1672.5
$ws.Range("M100").Value = -333843.34
$ws.Range("N100").Value = -2754.5

# Row 107
$ws.Range("H107").Value = 55556236
$ws.Range("I107").Value = 83333940
$ws.Range("J107").Value = 842.6667
$ws.Range("K107").Value = 250001820
$ws.Range("L107").Value = 2528.0001
$ws.Range("M107").Value = -249999900
$ws.Range("N107").Value = -6368.0001

